# Updates cryptos list (Price / Volume(1h) columns) to match the latest
# scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = "<new price text>"; E = "<new volume text>" }
# (only columns that actually changed are listed per row)
$changes = @{
    2  = @{ D = "27.236.68";  E = "  -0.11%  " }
    3  = @{ D = "1.902.76";   E = "  +0.19%  " }
    4  = @{                   E = "  -0.07%  " }
    5  = @{ D = "306.49";     E = "  -0.53%  " }
    6  = @{                   E = "  -0.04%  " }
    7  = @{ D = "0.5342";     E = "  +2.56%  " }
    8  = @{ D = "0.3808" }
    9  = @{ D = "0.07286";    E = "  +0.02%  " }
    10 = @{ D = "22.28";      E = "  +4.99%  " }
    11 = @{ D = "0.9025";     E = "  +0.15%  " }
    12 = @{ D = "0.08232";    E = "  +0.83%  " }
    13 = @{ D = "96.04";      E = "  -0.72%  " }
    14 = @{ D = "5.336";      E = "  +0.75%  " }
    15 = @{ D = "1.001";      E = "  -0.16%  " }
    16 = @{                   E = "  +1.89%  " }
    17 = @{                   E = "  +0.41%  " }
    18 = @{                   E = "  -0.04%  " }
    19 = @{ D = "27.261.65";  E = "  -0.11%  " }
    20 = @{ D = "5.033";      E = "  -1.20%  " }
    21 = @{ D = "1.096.24";   E = "  -42.36%  " }
    22 = @{ D = "10.77";      E = "  +0.54%  " }
    23 = @{ D = "6.502";      E = "  +1.35%  " }
    24 = @{ D = "149.63";     E = "  +1.57%  " }
    25 = @{ D = "2.298";      E = "  -0.22%  " }
    26 = @{ D = "18.36";      E = "  +0.66%  " }
    27 = @{                   E = "  +0.40%  " }
    28 = @{ D = "116.81";     E = "  +1.25%  " }
    29 = @{ D = "4.809";      E = "  -0.31%  " }
    30 = @{ D = "4.782";      E = "  -2.63%  " }
    31 = @{                   E = "  -0.17%  " }
    32 = @{ D = "0.8284";     E = "  +3.96%  " }
    33 = @{ D = "0.05063";    E = "  -0.10%  " }
    34 = @{ D = "1.220";      E = "  -0.91%  " }
    35 = @{ D = "2.998";      E = "  +1.14%  " }
    36 = @{                   E = "  -2.95%  " }
    37 = @{ D = "2.677";      E = "  +3.09%  " }
    38 = @{                   E = "  +0.98%  " }
    39 = @{ D = "0.02006";    E = "  +0.34%  " }
    40 = @{ D = "1.076";      E = "  -0.08%  " }
    41 = @{ D = "9.399";      E = "  +4.59%  " }
    42 = @{ D = "6.597";      E = "  +0.41%  " }
    43 = @{ D = "116.96";     E = "  +1.44%  " }
    44 = @{ D = "0.1523";     E = "  +0.43%  " }
    45 = @{ D = "0.4948";     E = "  +1.34%  " }
    47 = @{ D = "10.10";      E = "  +0.15%  " }
    48 = @{ D = "1.637";      E = "  +0.97%  " }
    49 = @{ D = "38.31";      E = "  +0.46%  " }
    50 = @{ D = "0.06180";    E = "  +3.87%  " }
    51 = @{ D = "63.32";      E = "  -0.68%  " }
}

foreach ($row in $changes.Keys) {
    $rowChanges = $changes[$row]
    foreach ($col in $rowChanges.Keys) {
        $cell = $ws.Range("$col$row")
        # Force the value to be written as text (matching the original
        # inline-string cells) instead of letting Excel auto-convert
        # numeric-looking strings (e.g. "306.49") into real numbers.
        # Clearing the formats afterwards removes the temporary "Text"
        # number-format style again so the cell keeps its original
        # (default) style, just like in the source workbook.
        $cell.NumberFormat = "@"
        $cell.Value = $rowChanges[$col]
        $cell.ClearFormats()
    }
}
